$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 1598.272
$ws.Range("C6").Value = 23922.944009123148
$ws.Range("C7").Value = 23369.110675789823
$ws.Range("C8").Value = 21388.59960821084
$ws.Range("C12").Value = 18810.61889584896
$ws.Range("C13").Value = 18414.61889584896
$ws.Range("C14").Value = 11682.618895848958
$ws.Range("C15").Value = 11183.075250348957
$ws.Range("C16").Value = 10867.11925034896
$ws.Range("C18").Value = 270.0
$ws.Range("C20").Value = 234603.93886706745
$ws.Range("C21").Value = 229172.68920873426
$ws.Range("C22").Value = 209750.5103478608
$ws.Range("C26").Value = 184469.15579497712
$ws.Range("C27").Value = 180585.7223949771
$ws.Range("C28").Value = 114567.35459497717
$ws.Range("C29").Value = 109668.50490383458
$ws.Range("C30").Value = 106570.03499643461
$ws.Range("C32").Value = 2647.7954999999993

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C2").Value = 5004.699999999999
$ws.Range("C6").Value = 1092.0
$ws.Range("D6").Value = -78.18051032029891
$ws.Range("C7").Value = 1365.0
$ws.Range("D7").Value = -72.72563790037364
$ws.Range("C8").Value = 2957.0
$ws.Range("D8").Value = -40.915539392970594
$ws.Range("C9").Value = 1131.0
$ws.Range("D9").Value = -77.40124283173816
$ws.Range("D10").Value = -50.226786820388824
$ws.Range("D11").Value = -26.109457110316285
$ws.Range("C12").Value = 2122.333333333333
$ws.Range("D12").Value = -57.593195729347734

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C2").Value = 3209.0999999999995
$ws.Range("C7").Value = 2346.0
$ws.Range("D7").Value = -26.8953912311863
$ws.Range("C8").Value = 1748.0
$ws.Range("D8").Value = -45.52989934872705
$ws.Range("C9").Value = 1048.0
$ws.Range("D9").Value = -67.34286871708578
$ws.Range("D10").Value = -28.35997631734753
$ws.Range("C11").Value = 2378.0
$ws.Range("D11").Value = -25.898226917204187
$ws.Range("C12").Value = 2184.0
$ws.Range("D12").Value = -31.94353557072075
$ws.Range("C13").Value = 1714.7142857142858
$ws.Range("D13").Value = -46.56712830032451

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C2").Value = 548.3
$ws.Range("D7").Value = -52.94546780959329
$ws.Range("D8").Value = -72.82509575050155
$ws.Range("C9").Value = 130.0
$ws.Range("D9").Value = -76.29035199708188
$ws.Range("C10").Value = 179.0
$ws.Range("D10").Value = -67.35363851905889

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C2").Value = 452.69999999999993
$ws.Range("D7").Value = -27.104042412193493
$ws.Range("C8").Value = 227.0
$ws.Range("D8").Value = -49.85641705323613
$ws.Range("C9").Value = 278.5
$ws.Range("D9").Value = -38.48022973271481

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C2").Value = 1306.7999999999997
$ws.Range("C3").Value = 451.9999999999999
$ws.Range("D3").Value = -88.47056422813996
$ws.Range("D9").Value = -64.0342822161004
$ws.Range("C10").Value = 249.0
$ws.Range("D10").Value = -61.891643709825516
$ws.Range("D11").Value = -70.30915212733393
$ws.Range("C12").Value = 226.0
$ws.Range("D16").Value = -64.0342822161004
$ws.Range("C17").Value = 249.0
$ws.Range("D17").Value = -61.891643709825516
$ws.Range("D18").Value = -70.30915212733393
$ws.Range("C19").Value = 226.0

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C2").Value = 2233.0999999999995
$ws.Range("C5").Value = 790.0
$ws.Range("D5").Value = -64.62316958488199
$ws.Range("C6").Value = 958.0
$ws.Range("D6").Value = -57.099995521920185
$ws.Range("C7").Value = 1086.0
$ws.Range("D7").Value = -51.368053378711195
$ws.Range("C8").Value = 941.0
$ws.Range("D8").Value = -57.86126908781513
$ws.Range("C9").Value = 943.75
$ws.Range("D9").Value = -57.73812189333212

$ws = $wb.Worksheets.Item("SYSTEMS")
$ws.Range("C2").Value = 1822.772
$ws.Range("D5").Value = 35.89192724048866
$ws.Range("D6").Value = 35.889446747116224
